$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 5

$ws.Cells.Item($row, 1).Value = 42606.881053240744
$ws.Cells.Item($row, 2).Value = 20
$ws.Cells.Item($row, 3).Value = 59
$ws.Cells.Item($row, 4).Value = 39
$ws.Cells.Item($row, 5).Value = 61
$ws.Cells.Item($row, 6).Value = 38
$ws.Cells.Item($row, 7).Value = 9082
$ws.Cells.Item($row, 8).Value = 6988
$ws.Cells.Item($row, 9).Value = 1203
$ws.Cells.Item($row, 10).Value = 115
$ws.Cells.Item($row, 11).Value = 76
$ws.Cells.Item($row, 12).Value = 8
$ws.Cells.Item($row, 13).Value = 5
$ws.Cells.Item($row, 14).Value = "Noun"

# Copy the date format (style index) from the cell above onto the new date cell,
# so it matches the existing date-formatted column without adding a new style.
$ws.Cells.Item($row - 1, 1).Copy() | Out-Null
$ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
